$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last row (old row 7) — the data now only spans rows 1-6.
$ws.Rows("7:7").Delete()

# New row data (TPM-updated values). Columns A-D are the categorical
# labels (Sending cluster, Ligand symbol, Receptor symbol, Target cluster);
# columns E-T are the numeric measures.
$rows = @(
  @{ Row=2;  A="MuSCs"; B="Fgf16"; C="Fgfr3"; D="ECs";
     E=3; F=1; G=0.2782933333333333; H=0.83488; I=1; J=1; K=3; L=1;
     M=4.662797333333334; N=13.988392; O=0.7324994586787992; P=0.7324994586787993;
     Q=1.297625412551111; R=11.67862871296; S=0.7324994586787992; T=0.7324994586787993 },
  @{ Row=3;  A="MuSCs"; B="Fgf16"; C="Fgfr3"; D="FAPs";
     E=3; F=1; G=0.2782933333333333; H=0.83488; I=1; J=1; K=3; L=1;
     M=0.6655859999999999; N=1.996758; O=0.1045598489170565; P=0.1045598489170565;
     Q=0.18522814656; R=1.66705331904; S=0.1045598489170565; T=0.1045598489170565 },
  @{ Row=4;  A="MuSCs"; B="Fgf16"; C="Fgfr3"; D="MuSCs";
     E=3; F=1; G=0.2782933333333333; H=0.83488; I=1; J=1; K=3; L=1;
     M=0.7894166666666665; N=2.36825; O=0.1240129561007488; P=0.1240129561007488;
     Q=0.2196893955555555; R=1.97720456; S=0.1240129561007488; T=0.1240129561007488 },
  @{ Row=5;  A="MuSCs"; B="Fgf16"; C="Fgfr3"; D="Neutrophils";
     E=3; F=1; G=0.2782933333333333; H=0.83488; I=1; J=1; K=1; L=0.3333333333333333;
     M=0.08057833333333334; N=0.241735; O=0.01265840681643176; P=0.01265840681643176;
     Q=0.02242441297777778; R=0.2018197168; S=0.01265840681643176; T=0.01265840681643176 },
  @{ Row=6;  A="MuSCs"; B="Fgf16"; C="Fgfr3"; D="Resolving-Mac";
     E=3; F=1; G=0.2782933333333333; H=0.83488; I=1; J=1; K=1; L=0.3333333333333333;
     M=0.16722; N=0.50166; O=0.02626932948696365; P=0.02626932948696365;
     Q=0.0465362112; R=0.4188259008; S=0.02626932948696365; T=0.02626932948696365 }
)

foreach ($r in $rows) {
  $row = $r.Row
  $ws.Cells.Item($row, 1).Value = $r.A
  $ws.Cells.Item($row, 2).Value = $r.B
  $ws.Cells.Item($row, 3).Value = $r.C
  $ws.Cells.Item($row, 4).Value = $r.D
  $ws.Cells.Item($row, 5).Value = $r.E
  $ws.Cells.Item($row, 6).Value = $r.F
  $ws.Cells.Item($row, 7).Value = $r.G
  $ws.Cells.Item($row, 8).Value = $r.H
  $ws.Cells.Item($row, 9).Value = $r.I
  $ws.Cells.Item($row, 10).Value = $r.J
  $ws.Cells.Item($row, 11).Value = $r.K
  $ws.Cells.Item($row, 12).Value = $r.L
  $ws.Cells.Item($row, 13).Value = $r.M
  $ws.Cells.Item($row, 14).Value = $r.N
  $ws.Cells.Item($row, 15).Value = $r.O
  $ws.Cells.Item($row, 16).Value = $r.P
  $ws.Cells.Item($row, 17).Value = $r.Q
  $ws.Cells.Item($row, 18).Value = $r.R
  $ws.Cells.Item($row, 19).Value = $r.S
  $ws.Cells.Item($row, 20).Value = $r.T
}
